# Update "snapshot" sheet: refresh the scraped_at timestamps (column K, rows 2-32)
$wb = $excel.ActiveWorkbook

$snapshot = $wb.Worksheets.Item("snapshot")

$newScrapedAt = @{
    2  = "2025-11-27T07:02:15.813256+00:00"
    3  = "2025-11-27T07:02:18.505254+00:00"
    4  = "2025-11-27T07:02:18.505287+00:00"
    5  = "2025-11-27T07:02:18.505307+00:00"
    6  = "2025-11-27T07:02:21.358081+00:00"
    7  = "2025-11-27T07:02:24.168007+00:00"
    8  = "2025-11-27T07:02:26.589564+00:00"
    9  = "2025-11-27T07:02:30.008443+00:00"
    10 = "2025-11-27T07:02:35.012012+00:00"
    11 = "2025-11-27T07:02:37.948007+00:00"
    12 = "2025-11-27T07:02:40.266757+00:00"
    13 = "2025-11-27T07:02:45.572430+00:00"
    14 = "2025-11-27T07:02:45.572462+00:00"
    15 = "2025-11-27T07:02:45.572487+00:00"
    16 = "2025-11-27T07:02:48.217316+00:00"
    17 = "2025-11-27T07:02:48.217351+00:00"
    18 = "2025-11-27T07:02:48.217372+00:00"
    19 = "2025-11-27T07:02:48.217391+00:00"
    20 = "2025-11-27T07:02:50.538544+00:00"
    21 = "2025-11-27T07:02:50.538576+00:00"
    22 = "2025-11-27T07:02:52.826512+00:00"
    23 = "2025-11-27T07:02:52.826543+00:00"
    24 = "2025-11-27T07:02:52.826562+00:00"
    25 = "2025-11-27T07:02:52.826581+00:00"
    26 = "2025-11-27T07:02:55.587073+00:00"
    27 = "2025-11-27T07:02:58.365887+00:00"
    28 = "2025-11-27T07:03:01.416339+00:00"
    29 = "2025-11-27T07:03:01.416374+00:00"
    30 = "2025-11-27T07:03:06.486007+00:00"
    31 = "2025-11-27T07:03:08.783400+00:00"
    32 = "2025-11-27T07:03:08.783429+00:00"
}

foreach ($row in $newScrapedAt.Keys) {
    $snapshot.Cells.Item($row, 11).Value = $newScrapedAt[$row]
}

# "new_injured" sheet: the previously-staged new-injury row got processed,
# so remove it, leaving just the header row.
$newInjured = $wb.Worksheets.Item("new_injured")
$newInjured.Rows.Item(2).Delete()
